$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window/view changes
$excel.ActiveWindow.WindowState = -4143  # xlNormal, ensure baseline

# Set formulas for each 10-row block (with shorter 9/8-row blocks near the end)
$blocks = @(
    @(1,10), @(11,20), @(21,30), @(31,40), @(41,50), @(51,60), @(61,69), @(70,78),
    @(79,87), @(88,96), @(97,105), @(106,115), @(116,125), @(126,135), @(136,145),
    @(146,155), @(156,165), @(166,175), @(176,183), @(184,191)
)

foreach ($blk in $blocks) {
    $s = $blk[0]
    $e = $blk[1]
    $startRow = $s + 2
    if ($startRow -le $e) {
        $rng = $ws.Range("B$startRow`:B$e")
        $rng.FormulaR1C1 = "=R[-1]C+0.85"
    }
}

# Sheet view changes
$ws.Range("K181").Select()
$excel.ActiveWindow.ScrollRow = 107

$wb.Save()
